$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header fields ---
$ws.Range("E11").Value2 = 321880
$ws.Range("C13").Value2 = 4

# --- Insert new row 21 for the 6th worker-period entry, matching the previous last-row style ---
# (this also shifts the old rows 25/26 down to 26/27)
$ws.Rows.Item(21).Insert()
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# --- Update worker-period table rows 16-21 ---
$ws.Range("E16").Value2 = "2506"
$ws.Range("G16").Value2 = 781242

$ws.Range("E17").Value2 = "2503"
$ws.Range("G17").Value2 = 781242

$ws.Range("E18").Value2 = "2502"
$ws.Range("G18").Value2 = 781242

$ws.Range("C19").Value2 = "1007655694"
$ws.Range("D19").Value2 = "MANOLO JOSE BABILONIA ACEVEDO"
$ws.Range("E19").Value2 = "2409"

$ws.Range("C20").Value2 = "1047510337"
$ws.Range("D20").Value2 = "ANGE MARCELA MONTES MATURANA"
$ws.Range("E20").Value2 = "2507"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

$ws.Range("C21").Value2 = "1143380075"
$ws.Range("D21").Value2 = "JHAXLYN NATALIA NARCISA ROCERO HERNANDEZ"
$ws.Range("E21").Value2 = "2507"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500

# --- Rework the signature block (now at rows 26/27 after the earlier row insert): ---
# remove old "___" line (row26), shift NOMBRE/FIRMA down to row27, add fresh "___" row26
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(27).Insert()

$ws.Range("B27:C27").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("H27:J27").Copy()
$ws.Range("H26:J26").PasteSpecial(-4122)
$ws.Range("B26:C26").Merge()
$ws.Range("H26:J26").Merge()

$ws.Range("B26").Value2 = "___________________________________"
$ws.Range("H26").Value2 = "___________________________________"
